$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values must stay as literal text (preserve trailing
# zeros / leading-zero decimal formatting), so force text number format
# before assigning, otherwise Excel auto-coerces numeric-looking strings
# into floating point numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.898.58'
$ws.Range('E2').Value = '  +1.15%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.627.93'
$ws.Range('E3').Value = '  +1.93%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.13'
$ws.Range('E5').Value = '  +0.94%  '
$ws.Range('E6').Value = '  +1.06%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '29.71'
$ws.Range('E8').Value = '  +10.53%  '
$ws.Range('E9').Value = '  +3.34%  '
$ws.Range('E10').Value = '  +2.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0916'
$ws.Range('E11').Value = '  +0.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.861.07'
$ws.Range('E12').Value = '  +1.99%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.639.05'
$ws.Range('E13').Value = '  +2.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.570'
$ws.Range('E14').Value = '  +6.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '9.17'
$ws.Range('E15').Value = '  +20.95%  '
$ws.Range('E16').Value = '  +4.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.903.27'
$ws.Range('E17').Value = '  +1.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '64.96'
$ws.Range('E18').Value = '  +1.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '248.65'
$ws.Range('E19').Value = '  +2.91%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0704'
$ws.Range('E20').Value = '  +1.78%  '
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('E22').Value = '  +4.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.60'
$ws.Range('E23').Value = '  +4.17%  '
$ws.Range('E24').Value = '  +0.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '159.11'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.70'
$ws.Range('E26').Value = '  +2.10%  '
$ws.Range('E27').Value = '  +2.31%  '
$ws.Range('E28').Value = '  +3.17%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('E30').Value = '  +2.87%  '
$ws.Range('E31').Value = '  +5.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.36'
$ws.Range('E32').Value = '  +4.68%  '
$ws.Range('E33').Value = '  +1.77%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.427.86'
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.65'
$ws.Range('E35').Value = '  +7.33%  '
$ws.Range('E36').Value = '  +0.79%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.88'
$ws.Range('E37').Value = '  +0.24%  '
$ws.Range('E38').Value = '  -0.38%  '
$ws.Range('E39').Value = '  +3.14%  '
$ws.Range('E40').Value = '  +2.64%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '71.56'
$ws.Range('E41').Value = '  +9.25%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0500'
$ws.Range('E42').Value = '  +1.61%  '
$ws.Range('B43').Value = 'BitcoinSV'
$ws.Range('C43').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '55.23'
$ws.Range('E43').Value = '  +1.03%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.830'
$ws.Range('E44').Value = '  +3.81%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.97'
$ws.Range('E45').Value = '  +0.67%  '
$ws.Range('E46').Value = '  +5.68%  '
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.767.16'
$ws.Range('E49').Value = '  +1.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '89.61'
$ws.Range('E50').Value = '  +4.29%  '
$ws.Range('E51').Value = '  +10.64%  '
